$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Run mode for Notification module (row 7) from Y to N
$ws.Range("C7").Value = "N"

# Select the updated cell, matching the selection state after the edit
$ws.Activate()
$ws.Range("C7").Select()
